$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / info block edits ---
$ws.Range("A1").Value = "fimra01"

$ws.Range("B2").Value = "22.4.2020"
$ws.Range("D2").Value = "Nr faktury: "
$ws.Range("E2").Value = "assdff123"

$ws.Range("B3").Value = "30.4.2020"
$ws.Range("D3").Value = "Płatność: "
$ws.Range("E3").Value = 123

# --- Row 5 edits ---
$ws.Range("B5").Value = 65
$ws.Range("C5").Value = 21

# --- Row 6 edit (quantity) ---
$ws.Range("C6").Value = 32

# Break existing merges before reshaping the M1 group (rows 5-6 -> 5-7)
$ws.Range("A5:A6").UnMerge()
$ws.Range("B5:B6").UnMerge()

# Row 7 (new row, completes the M1 group)
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C7").Value = 12
$ws.Range("D7").Value = "bialy"

# Re-merge the M1 group across rows 5-7
$ws.Range("A5:A7").Merge()
$ws.Range("B5:B7").Merge()

# --- Rows 8-9: M3 group ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A8").Value = "M3"
$ws.Range("B8").Value = 44
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = "styropian"

$ws.Range("A6").Copy() | Out-Null
$ws.Range("A9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C9").Value = 32
$ws.Range("D9").Value = "czarny"

$ws.Range("A8:A9").Merge()
$ws.Range("B8:B9").Merge()

# --- Rows 10-11: Statyw drewniany group ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A10").Value = "Statyw drewniany"
$ws.Range("B10").Value = 55
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = "biały"

$ws.Range("A6").Copy() | Out-Null
$ws.Range("A11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C11").Value = 22
$ws.Range("D11").Value = "czarny"

$ws.Range("A10:A11").Merge()
$ws.Range("B10:B11").Merge()

# --- Rows 12-13: Statyw metalowy group ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A12").Value = "Statyw metalowy"
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 12
$ws.Range("D12").Formula = '="45"'

$ws.Range("A6").Copy() | Out-Null
$ws.Range("A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C13").Value = 3
$ws.Range("D13").Formula = '="90"'

$ws.Range("A12:A13").Merge()
$ws.Range("B12:B13").Merge()

$excel.CutCopyMode = 0
